$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "56.286.98"
$ws.Range("E2").Value = "  +2.00%  "
$ws.Range("D3").Value = "2.309.34"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue "D5" "516.51"
$ws.Range("E5").Value = "  +2.09%  "
Set-TextValue "D6" "134.33"
$ws.Range("E6").Value = "  +3.87%  "
Set-TextValue "D7" "0.997"
$ws.Range("E7").Value = "  +0.17%  "
Set-TextValue "D8" "0.535"
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("D9").Value = "2.330.46"
$ws.Range("E9").Value = "  +0.86%  "
Set-TextValue "D10" "0.102"
$ws.Range("E10").Value = "  +3.84%  "
Set-TextValue "D11" "0.152"
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("E12").Value = "  +3.97%  "
Set-TextValue "D13" "0.340"
$ws.Range("E13").Value = "  -0.54%  "
Set-TextValue "D14" "23.80"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").Value = "2.726.45"
$ws.Range("D16").Value = "56.453.78"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").Value = "2.334.07"
$ws.Range("E18").Value = "  +1.64%  "
Set-TextValue "D19" "10.43"
$ws.Range("E19").Value = "  -0.08%  "
Set-TextValue "D20" "4.21"
$ws.Range("E20").Value = "  +0.49%  "
Set-TextValue "D21" "321.42"
$ws.Range("E21").Value = "  +2.81%  "
Set-TextValue "D22" "6.53"
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("E23").Value = "  +0.27%  "
Set-TextValue "D24" "60.46"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("E25").Value = "  +4.83%  "
Set-TextValue "D26" "0.991"
$ws.Range("E26").Value = "  -0.23%  "
Set-TextValue "D27" "8.00"
$ws.Range("E27").Value = "  +6.25%  "
Set-TextValue "D28" "1.28"
$ws.Range("E28").Value = "  +10.48%  "
$ws.Range("D29").Value = "0.0₃0736"
$ws.Range("E29").Value = "  +3.45%  "
$ws.Range("E30").Value = "  +3.67%  "
Set-TextValue "D31" "166.60"
$ws.Range("E31").Value = "  -3.27%  "
Set-TextValue "D32" "6.18"
$ws.Range("E32").Value = "  +0.64%  "
Set-TextValue "D33" "18.30"
$ws.Range("E33").Value = "  +1.60%  "
Set-TextValue "D35" "0.995"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("E36").Value = "  +0.87%  "
Set-TextValue "D37" "0.915"
$ws.Range("E37").Value = "  +0.33%  "
Set-TextValue "D38" "4.01"
$ws.Range("E38").Value = "  +2.78%  "
Set-TextValue "D39" "1.54"
$ws.Range("E39").Value = "  +5.45%  "
Set-TextValue "D40" "37.90"
$ws.Range("E40").Value = "  +2.63%  "
Set-TextValue "D41" "0.380"
$ws.Range("E41").Value = "  +0.98%  "
Set-TextValue "D42" "139.08"
$ws.Range("E42").Value = "  +1.86%  "
Set-TextValue "D43" "3.59"
$ws.Range("E43").Value = "  +2.87%  "
$ws.Range("E44").Value = "  +3.76%  "
Set-TextValue "D45" "276.74"
$ws.Range("E45").Value = "  +6.08%  "
Set-TextValue "D46" "0.0930"
$ws.Range("E46").Value = "  +1.42%  "
Set-TextValue "D47" "0.0505"
$ws.Range("E47").Value = "  -0.45%  "
Set-TextValue "D48" "0.559"
$ws.Range("E48").Value = "  +1.45%  "
Set-TextValue "D49" "0.0216"
$ws.Range("E49").Value = "  +1.94%  "
Set-TextValue "D50" "0.379"
$ws.Range("E50").Value = "  +1.11%  "
Set-TextValue "D51" "17.76"
$ws.Range("E51").Value = "  +7.89%  "
